# Applies the commit's content edits to the workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"   # Date
$meta.Range("B15").Value = "4.0.1"                       # FHIR Version

# --- Elements sheet ---
$elems = $wb.Worksheets.Item("Elements")

# Row 2 (Extension) - Constraint(s): drop the "unless an empty Parameters resource ... or $this is Parameters" clause
$elems.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 (Extension.id) - Type(s): "id" -> "string"
$elems.Range("K3").Value = "string`n"

# Row 6 (Extension.value[x]) - Definition: R4B -> R4 in the Extensibility link
$elems.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
